{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Replaces the text of 7 paragraphs in the \"Code of Conduct\" document with\n// updated wording, per the commit \"Updated my part of Code of Conduct\".\n// Each target paragraph consists of a single run, so we match paragraphs by\n// their (unique) original prefix and overwrite the paragraph's text in\n// place with insertText(..., Word.InsertLocation.replace). This preserves\n// paragraph/run formatting (pPr / rPr) while swapping only the visible text.\n\nconst replacements = [\n  {\n    find: \"Results Based on Seriousness: The repercussions for not following the rules are contingent upon the seriousness of the transgression. Essentially, the severity of the consequences aligns with the gravity of the issue at hand. This approach ensures a fair and proportional response to rule violations, emphasizing the importance of context in determining appropriate actions.\",\n    replace: \"Results based on the severity: The repercussions for not following the rules correspond to how bad the action is. We are taking this approach to ensure a fair and proportional response to any violation of the rules. The severity of a mistake will be decided by the entirety of the team who are not part of the issue itself.\"\n  },\n  {\n    find: \"Small Mistakes: In instances of minor errors, the corrective action involves issuing a warning. This serves as a constructive measure aimed at fostering a learning environment. Individuals who make small mistakes receive guidance to help them understand and rectify their errors, facilitating continuous improvement and skill development.\",\n    replace: \"Small Mistakes: In instances of minor errors, we will issue a warning. These warnings have the purpose of inspiring corrective actions for those who make small mistakes. Although one warning may not be significant, not correcting it may generate harsher consequences.\"\n  },\n  {\n    find: \"Moderate Mistakes: For medium-sized mistakes, a structured approach is implemented \\u2013 the three-strikes rule. If an individual repeats the same mistake three times, a thorough evaluation ensues. This evaluation involves a discussion on whether the person should continue participating in the project. The objective is to ensure a cohesive and functional team, where repeated medium mistakes prompt a thoughtful reconsideration of a team member's fit within the project.\",\n    replace: \"Moderate Mistakes: For medium-sized mistakes, a three-strikes rule shall be implemented. If an individual reaches three strikes, a thorough evaluation will take place, which may very well lead to excluding that person from the project entirely. Given that there are three strikes, everyone will have a chance to correct their mistake and thus remove the corresponding strike.\"\n  },\n  {\n    find: \"Big Mistakes: Significant errors with the potential to cause substantial problems are treated with zero tolerance. In the case of substantial blunders, serious actions are taken. This firm stance is adopted because significant mistakes can have profound implications for the team. Swift and decisive measures are necessary to mitigate harm and uphold the overall integrity and progress of the project.\",\n    replace: \"Big Mistakes: Significant errors will be met with zero tolerance, as they can make or break the effort and dynamic of the entire group. These types of issues have to be dealt with quickly, to ensure regaining stability and progress in the project. There will be no chance to correct these actions, as their existence is not without malice or intent. Removal from the project will take place immediately\"\n  },\n  {\n    find: \"Before Meetings: Before our required meetings, team members talk separately to plan what they want to discuss. This helps us be ready and organized for the main meeting.\",\n    replace: \"Before Meetings: Before our required meetings, team members talk separately to plan what they want to discuss. This helps us be ready for anything regarding the agenda and organized for the main meeting, ensuring everyone is contributing.\"\n  },\n  {\n    find: \"Weekly Coding Talks: We have regular meetings to discuss our coding work every week. We talk about the tasks we need to do, who is responsible for what (Assignees), and other important things. This helps us keep track of our progress.\",\n    replace: \"Weekly Coding Talks: We have regular meetings to discuss our coding work every week. We talk about the tasks we need to do, who is responsible for what (Assignees), and other important things. This helps us keep track of our progress and ensures the contribution is split equally. (both in terms of the amount of code, as well as front-end/back-end)\"\n  },\n  {\n    find: \"Urgent Meetings: If something really urgent comes up, we have a plan for quick meetings to solve the problem. This helps us handle important issues right away and keeps our work moving smoothly.\",\n    replace: \"Urgent Meetings: If something really urgent comes up, we have a plan for quick meetings to solve the problem. This helps us handle important issues right away and keeps our workflow smooth.\"\n  }\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (const item of replacements) {\n  let matched = false;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === item.find) {\n      paragraphs.items[i].insertText(item.replace, Word.InsertLocation.replace);\n      matched = true;\n      break;\n    }\n  }\n  if (!matched) {\n    throw new Error(\"Paragraph to replace not found: \" + item.find.substring(0, 50));\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Replaces the text of 7 paragraphs in the \"Code of Conduct\" document with\n# updated wording, per the commit \"Updated my part of Code of Conduct\".\n# Each target paragraph is a single run; we match paragraphs by their\n# (unique) original text and overwrite the paragraph Range's Text in place.\n# This preserves paragraph/run formatting (pPr / rPr) while swapping only\n# the visible text, and leaves the trailing paragraph mark untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{\n    Find = \"Results Based on Seriousness: The repercussions for not following the rules are contingent upon the seriousness of the transgression. Essentially, the severity of the consequences aligns with the gravity of the issue at hand. This approach ensures a fair and proportional response to rule violations, emphasizing the importance of context in determining appropriate actions.\"\n    Replace = \"Results based on the severity: The repercussions for not following the rules correspond to how bad the action is. We are taking this approach to ensure a fair and proportional response to any violation of the rules. The severity of a mistake will be decided by the entirety of the team who are not part of the issue itself.\"\n  },\n  @{\n    Find = \"Small Mistakes: In instances of minor errors, the corrective action involves issuing a warning. This serves as a constructive measure aimed at fostering a learning environment. Individuals who make small mistakes receive guidance to help them understand and rectify their errors, facilitating continuous improvement and skill development.\"\n    Replace = \"Small Mistakes: In instances of minor errors, we will issue a warning. These warnings have the purpose of inspiring corrective actions for those who make small mistakes. Although one warning may not be significant, not correcting it may generate harsher consequences.\"\n  },\n  @{\n    Find = \"Moderate Mistakes: For medium-sized mistakes, a structured approach is implemented \u2013 the three-strikes rule. If an individual repeats the same mistake three times, a thorough evaluation ensues. This evaluation involves a discussion on whether the person should continue participating in the project. The objective is to ensure a cohesive and functional team, where repeated medium mistakes prompt a thoughtful reconsideration of a team member's fit within the project.\"\n    Replace = \"Moderate Mistakes: For medium-sized mistakes, a three-strikes rule shall be implemented. If an individual reaches three strikes, a thorough evaluation will take place, which may very well lead to excluding that person from the project entirely. Given that there are three strikes, everyone will have a chance to correct their mistake and thus remove the corresponding strike.\"\n  },\n  @{\n    Find = \"Big Mistakes: Significant errors with the potential to cause substantial problems are treated with zero tolerance. In the case of substantial blunders, serious actions are taken. This firm stance is adopted because significant mistakes can have profound implications for the team. Swift and decisive measures are necessary to mitigate harm and uphold the overall integrity and progress of the project.\"\n    Replace = \"Big Mistakes: Significant errors will be met with zero tolerance, as they can make or break the effort and dynamic of the entire group. These types of issues have to be dealt with quickly, to ensure regaining stability and progress in the project. There will be no chance to correct these actions, as their existence is not without malice or intent. Removal from the project will take place immediately\"\n  },\n  @{\n    Find = \"Before Meetings: Before our required meetings, team members talk separately to plan what they want to discuss. This helps us be ready and organized for the main meeting.\"\n    Replace = \"Before Meetings: Before our required meetings, team members talk separately to plan what they want to discuss. This helps us be ready for anything regarding the agenda and organized for the main meeting, ensuring everyone is contributing.\"\n  },\n  @{\n    Find = \"Weekly Coding Talks: We have regular meetings to discuss our coding work every week. We talk about the tasks we need to do, who is responsible for what (Assignees), and other important things. This helps us keep track of our progress.\"\n    Replace = \"Weekly Coding Talks: We have regular meetings to discuss our coding work every week. We talk about the tasks we need to do, who is responsible for what (Assignees), and other important things. This helps us keep track of our progress and ensures the contribution is split equally. (both in terms of the amount of code, as well as front-end/back-end)\"\n  },\n  @{\n    Find = \"Urgent Meetings: If something really urgent comes up, we have a plan for quick meetings to solve the problem. This helps us handle important issues right away and keeps our work moving smoothly.\"\n    Replace = \"Urgent Meetings: If something really urgent comes up, we have a plan for quick meetings to solve the problem. This helps us handle important issues right away and keeps our workflow smooth.\"\n  }\n)\n\nforeach ($pair in $pairs) {\n  $matched = $false\n  foreach ($p in $d.Paragraphs) {\n    $current = $p.Range.Text.TrimEnd([char]13)\n    if ($current -eq $pair.Find) {\n      $p.Range.Text = $pair.Replace\n      $matched = $true\n      break\n    }\n  }\n  if (-not $matched) {\n    throw \"Paragraph to replace not found: $($pair.Find.Substring(0, 40))\"\n  }\n}\n"}
